$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("A19").Copy()
$ws1.Range("A20").PasteSpecial(-4122)
$ws1.Range("A20").Value = 18
$ws1.Range("B2").Value = '''2024-05-18'
$ws1.Range("C2").Value = '太仓·原x崩铁ONLY'
$ws1.Range("D2").Value = '滨河路128号 凯景世纪大酒店(太仓滨河路店)'
$ws1.Range("E2").Value = '2024.05.18 10:00-05.18 17:00'
$ws1.Range("F2").Value = 59
$ws1.Range("G2").Value = 55
$ws1.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=84904'
$ws1.Range("I2").Value = '//i2.hdslb.com/bfs/openplatform/202404/qxpmdgLs1714118849736.jpeg'
$ws1.Range("B3").Value = '''2024-05-18'
$ws1.Range("C3").Value = '苏州·OrangeOrange国潮&随机宅舞派对【免费活动】'
$ws1.Range("D3").Value = '狮山路298号 金鹰国际购物中心(狮山路店)'
$ws1.Range("E3").Value = '2024.05.18 13:00-05.18 17:00'
$ws1.Range("F3").Value = 118
$ws1.Range("G3").Value = 29
$ws1.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=83949'
$ws1.Range("I3").Value = '//i1.hdslb.com/bfs/openplatform/202404/DOH6BK8i1712638105049.png'
$ws1.Range("B4").Value = '''2024-05-18'
$ws1.Range("C4").Value = '苏州·YoungComic动漫嘉年华'
$ws1.Range("D4").Value = '清禾路886号 尹山湖大剧院'
$ws1.Range("E4").Value = '2024.05.18 10:00-05.18 17:00'
$ws1.Range("F4").Value = 1733
$ws1.Range("G4").Value = 60
$ws1.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=83142'
$ws1.Range("I4").Value = '//i2.hdslb.com/bfs/openplatform/202403/4wWLK6Jg1710840463319.jpeg'
$ws1.Range("B5").Value = '''2024-05-19'
$ws1.Range("C5").Value = '苏州·国乙only（免费展）'
$ws1.Range("D5").Value = '吴中万达广场 吴中万达广场'
$ws1.Range("E5").Value = '2024.05.19 14:00-05.19 17:00'
$ws1.Range("F5").Value = 303
$ws1.Range("G5").Value = 20
$ws1.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=85022'
$ws1.Range("I5").Value = '//i0.hdslb.com/bfs/openplatform/202404/4kJ5GWEo1714137407259.jpeg'
$ws1.Range("B6").Value = '''2024-05-19'
$ws1.Range("C6").Value = '苏州·国潮动漫聚会'
$ws1.Range("D6").Value = '金门路33号(金门路阊胥路交界) 长船湾青年码头运河剧场'
$ws1.Range("E6").Value = '2024.05.19 09:00-05.19 18:00'
$ws1.Range("F6").Value = 13
$ws1.Range("G6").Value = 20
$ws1.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=85372'
$ws1.Range("I6").Value = '//i0.hdslb.com/bfs/openplatform/202405/KX0e5kR71715239777929.jpeg'
$ws1.Range("B7").Value = '''2024-05-25'
$ws1.Range("C7").Value = '苏州·姑苏梦行高校联展'
$ws1.Range("D7").Value = '尹山湖商业水街2号楼3层 格莱美婚礼宴会中心(尹山湖商业水街店)'
$ws1.Range("E7").Value = '2024.05.25 10:00-05.25 17:30'
$ws1.Range("F7").Value = 81
$ws1.Range("G7").Value = 60
$ws1.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=84991'
$ws1.Range("I7").Value = '//i0.hdslb.com/bfs/openplatform/202404/ZSJS1FTx1713888009395.png'
$ws1.Range("B8").Value = '''2024-05-25'
$ws1.Range("C8").Value = '苏州·燃梦Project'
$ws1.Range("D8").Value = '清禾路886号 尹山湖大剧院'
$ws1.Range("E8").Value = '2024.05.25 10:30-05.25 16:30'
$ws1.Range("F8").Value = 2039
$ws1.Range("G8").Value = 60
$ws1.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=83926'
$ws1.Range("I8").Value = '//i2.hdslb.com/bfs/openplatform/202404/T3neM9fF1714119047940.jpeg'
$ws1.Range("B9").Value = '''2024-06-08'
$ws1.Range("C9").Value = '【会员购严选】苏州·Come in joy动漫国潮文化节'
$ws1.Range("D9").Value = '金山南路288号 广电国际会展中心'
$ws1.Range("E9").Value = '2024.06.08 10:00-06.09 17:00'
$ws1.Range("F9").Value = 10311
$ws1.Range("G9").Value = 60
$ws1.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=82233'
$ws1.Range("I9").Value = '//i0.hdslb.com/bfs/openplatform/202403/F86lgbSt1709278264141.jpeg'
$ws1.Range("B10").Value = '''2024-06-15'
$ws1.Range("C10").Value = '苏州·蔚蓝档案ONLY#2024~Game Builders Go!!!!'
$ws1.Range("D10").Value = '城际路21号 苏州汇融广场假日酒店'
$ws1.Range("E10").Value = '2024.06.15 10:00-06.15 17:00'
$ws1.Range("F10").Value = 176
$ws1.Range("G10").Value = 75
$ws1.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=84130'
$ws1.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202404/bpTzFcDq1713253785881.jpeg'
$ws1.Range("B11").Value = '''2024-06-16'
$ws1.Range("C11").Value = '苏州·明日方舟ONLY#2024~佑桑柔'
$ws1.Range("D11").Value = '城际路21号 苏州汇融广场假日酒店'
$ws1.Range("E11").Value = '2024.06.16 10:00-06.16 17:00'
$ws1.Range("F11").Value = 147
$ws1.Range("G11").Value = 75
$ws1.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=84046'
$ws1.Range("I11").Value = '//i1.hdslb.com/bfs/openplatform/202405/FxnSAnJm1714978943173.jpeg'
$ws1.Range("B12").Value = '''2024-06-29'
$ws1.Range("C12").Value = '苏州·归离之缘原神only展'
$ws1.Range("D12").Value = '清禾路888号2号楼3楼 格莱美婚礼宴会中心'
$ws1.Range("E12").Value = '2024.06.29 09:30-06.29 18:30'
$ws1.Range("F12").Value = 265
$ws1.Range("G12").Value = 89
$ws1.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=83271'
$ws1.Range("I12").Value = '//i1.hdslb.com/bfs/openplatform/202403/hNkSoRCt1710999968954.png'
$ws1.Range("B13").Value = '''2024-07-06'
$ws1.Range("C13").Value = '苏州·第一届寒假动漫展宅舞比赛-CF01'
$ws1.Range("D13").Value = '润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店'
$ws1.Range("E13").Value = '2024.07.06 10:00-07.06 16:00'
$ws1.Range("F13").Value = 196
$ws1.Range("G13").Value = 49
$ws1.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=80528'
$ws1.Range("I13").Value = '//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg'
$ws1.Range("B14").Value = '''2024-07-20'
$ws1.Range("C14").Value = '苏州·白日梦想7.20全职猎人ONLY展'
$ws1.Range("D14").Value = '金芳路与新发路交叉口东南120米 万龙大厦'
$ws1.Range("E14").Value = '2024.07.20 09:00-07.20 17:00'
$ws1.Range("F14").Value = 394
$ws1.Range("G14").Value = 72
$ws1.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=83508'
$ws1.Range("I14").Value = '//i1.hdslb.com/bfs/openplatform/202403/LC3LtiWw1711633827120.jpeg'
$ws1.Range("B15").Value = '''2024-07-20'
$ws1.Range("C15").Value = '苏州·萤火国潮文化节动漫品牌博览会'
$ws1.Range("D15").Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws1.Range("E15").Value = '2024.07.20 10:00-07.21 17:00'
$ws1.Range("F15").Value = 7170
$ws1.Range("G15").Value = 60
$ws1.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=83301'
$ws1.Range("I15").Value = '//i0.hdslb.com/bfs/openplatform/202403/rV07luU61711274774556.jpeg'
$ws1.Range("B16").Value = '''2024-07-27'
$ws1.Range("C16").Value = '苏州·第一届动漫游戏展'
$ws1.Range("D16").Value = '清禾路886号 尹山湖大剧院'
$ws1.Range("E16").Value = '2024.07.27 10:30-07.27 17:00'
$ws1.Range("F16").Value = 1106
$ws1.Range("G16").Value = 60
$ws1.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=84899'
$ws1.Range("I16").Value = '//i2.hdslb.com/bfs/openplatform/202404/ARz0BVLv1712661597595.jpeg'
$ws1.Range("B17").Value = '''2024-08-03'
$ws1.Range("C17").Value = '苏州·星部落动漫嘉年华'
$ws1.Range("D17").Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Range("E17").Value = '2024.08.03 09:00-08.04 16:00'
$ws1.Range("F17").Value = 676
$ws1.Range("G17").Value = 49
$ws1.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=84858'
$ws1.Range("I17").Value = '//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg'
$ws1.Range("B18").Value = '''2024-08-17'
$ws1.Range("C18").Value = '苏州·ICAN summer World动漫品牌夏游节'
$ws1.Range("D18").Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws1.Range("E18").Value = '2024.08.17 10:00-08.18 17:00'
$ws1.Range("F18").Value = 108
$ws1.Range("G18").Value = 60
$ws1.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=85289'
$ws1.Range("I18").Value = '//i0.hdslb.com/bfs/openplatform/202404/JavlW9fj1714459472747.jpeg'
$ws1.Range("B19").Value = '''2024-10-01'
$ws1.Range("C19").Value = '苏州·第二届Redamancy动漫游戏嘉年华'
$ws1.Range("D19").Value = '长江路436号绿宝广场b1 party king运动街区'
$ws1.Range("E19").Value = '2024.10.01 10:00-10.05 17:00'
$ws1.Range("F19").Value = 58
$ws1.Range("G19").Value = 98
$ws1.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=83576'
$ws1.Range("I19").Value = '//i2.hdslb.com/bfs/openplatform/202403/MKyrtd4c1711689984512.jpeg'
$ws1.Range("B20").Value = '''2024-10-01'
$ws1.Range("C20").Value = '苏州·第十三届理想乡动漫展-同人创作者大会'
$ws1.Range("D20").Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Range("E20").Value = '2024.10.01 10:00-10.03 17:00'
$ws1.Range("F20").Value = 254
$ws1.Range("G20").Value = 39
$ws1.Range("H20").Value = 'https://show.bilibili.com/platform/detail.html?id=83821'
$ws1.Range("I20").Value = '//i0.hdslb.com/bfs/openplatform/202404/OMtuTTFY1711958198579.jpeg'

# ---- Sheet: 演出 ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14

# ---- Sheet: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("A22").Copy()
$ws4.Range("A23").PasteSpecial(-4122)
$ws4.Range("A23").Value = 21
$ws4.Range("B2").Value = '''2024-05-18'
$ws4.Range("C2").Value = '太仓·原x崩铁ONLY'
$ws4.Range("D2").Value = '滨河路128号 凯景世纪大酒店(太仓滨河路店)'
$ws4.Range("E2").Value = '2024.05.18 10:00-05.18 17:00'
$ws4.Range("F2").Value = 59
$ws4.Range("G2").Value = 55
$ws4.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=84904'
$ws4.Range("I2").Value = '//i2.hdslb.com/bfs/openplatform/202404/qxpmdgLs1714118849736.jpeg'
$ws4.Range("B3").Value = '''2024-05-18'
$ws4.Range("C3").Value = '苏州·OrangeOrange国潮&随机宅舞派对【免费活动】'
$ws4.Range("D3").Value = '狮山路298号 金鹰国际购物中心(狮山路店)'
$ws4.Range("E3").Value = '2024.05.18 13:00-05.18 17:00'
$ws4.Range("F3").Value = 118
$ws4.Range("G3").Value = 29
$ws4.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=83949'
$ws4.Range("I3").Value = '//i1.hdslb.com/bfs/openplatform/202404/DOH6BK8i1712638105049.png'
$ws4.Range("B4").Value = '''2024-05-18'
$ws4.Range("C4").Value = '苏州·YoungComic动漫嘉年华'
$ws4.Range("D4").Value = '清禾路886号 尹山湖大剧院'
$ws4.Range("E4").Value = '2024.05.18 10:00-05.18 17:00'
$ws4.Range("F4").Value = 1733
$ws4.Range("G4").Value = 60
$ws4.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=83142'
$ws4.Range("I4").Value = '//i2.hdslb.com/bfs/openplatform/202403/4wWLK6Jg1710840463319.jpeg'
$ws4.Range("B5").Value = '''2024-05-19'
$ws4.Range("C5").Value = '苏州·国乙only（免费展）'
$ws4.Range("D5").Value = '吴中万达广场 吴中万达广场'
$ws4.Range("E5").Value = '2024.05.19 14:00-05.19 17:00'
$ws4.Range("F5").Value = 303
$ws4.Range("G5").Value = 20
$ws4.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=85022'
$ws4.Range("I5").Value = '//i0.hdslb.com/bfs/openplatform/202404/4kJ5GWEo1714137407259.jpeg'
$ws4.Range("B6").Value = '''2024-05-19'
$ws4.Range("C6").Value = '苏州·国潮动漫聚会'
$ws4.Range("D6").Value = '金门路33号(金门路阊胥路交界) 长船湾青年码头运河剧场'
$ws4.Range("E6").Value = '2024.05.19 09:00-05.19 18:00'
$ws4.Range("F6").Value = 14
$ws4.Range("G6").Value = 20
$ws4.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=85372'
$ws4.Range("I6").Value = '//i0.hdslb.com/bfs/openplatform/202405/KX0e5kR71715239777929.jpeg'
$ws4.Range("B7").Value = '''2024-05-19'
$ws4.Range("C7").Value = '苏州·跨越二次元ACG神级动漫世界巡回演唱会'
$ws4.Range("D7").Value = '吴中区清禾路886号 苏州聚橙尹山湖大剧院'
$ws4.Range("E7").Value = '2024.05.19 19:30-05.19 21:10'
$ws4.Range("F7").Value = 14
$ws4.Range("G7").Value = 60
$ws4.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=85078'
$ws4.Range("I7").Value = '//i1.hdslb.com/bfs/openplatform/202404/l1Y850En1714465446331.jpeg'
$ws4.Range("B8").Value = '''2024-05-25'
$ws4.Range("C8").Value = '苏州·姑苏梦行高校联展'
$ws4.Range("D8").Value = '尹山湖商业水街2号楼3层 格莱美婚礼宴会中心(尹山湖商业水街店)'
$ws4.Range("E8").Value = '2024.05.25 10:00-05.25 17:30'
$ws4.Range("F8").Value = 81
$ws4.Range("G8").Value = 60
$ws4.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=84991'
$ws4.Range("I8").Value = '//i0.hdslb.com/bfs/openplatform/202404/ZSJS1FTx1713888009395.png'
$ws4.Range("B9").Value = '''2024-05-25'
$ws4.Range("C9").Value = '苏州·燃梦Project'
$ws4.Range("D9").Value = '清禾路886号 尹山湖大剧院'
$ws4.Range("E9").Value = '2024.05.25 10:30-05.25 16:30'
$ws4.Range("F9").Value = 2039
$ws4.Range("G9").Value = 60
$ws4.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=83926'
$ws4.Range("I9").Value = '//i2.hdslb.com/bfs/openplatform/202404/T3neM9fF1714119047940.jpeg'
$ws4.Range("B10").Value = '''2024-06-01'
$ws4.Range("C10").Value = '苏州·春日计划2024——特别二次元不插电音乐会'
$ws4.Range("D10").Value = '星湖街555号高教区(体育馆南侧) 苏州独墅湖影剧院'
$ws4.Range("E10").Value = '2024.06.01 19:30-06.01 21:00'
$ws4.Range("F10").Value = 557
$ws4.Range("G10").Value = 88
$ws4.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=84720'
$ws4.Range("I10").Value = '//i1.hdslb.com/bfs/openplatform/202404/gwLWvSew1713796405109.png'
$ws4.Range("B11").Value = '''2024-06-02'
$ws4.Range("C11").Value = '苏州·英雄时代2024哈瓦西钢琴演奏会'
$ws4.Range("D11").Value = '东太湖大道12000号 苏州湾大剧院'
$ws4.Range("E11").Value = '2024.06.02 19:30-06.02 21:00'
$ws4.Range("F11").Value = 1
$ws4.Range("G11").Value = 499
$ws4.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=83901'
$ws4.Range("I11").Value = '//i0.hdslb.com/bfs/openplatform/202404/LbCirky11712569675168.png'
$ws4.Range("B12").Value = '''2024-06-08'
$ws4.Range("C12").Value = '【会员购严选】苏州·Come in joy动漫国潮文化节'
$ws4.Range("D12").Value = '金山南路288号 广电国际会展中心'
$ws4.Range("E12").Value = '2024.06.08 10:00-06.09 17:00'
$ws4.Range("F12").Value = 10311
$ws4.Range("G12").Value = 60
$ws4.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=82233'
$ws4.Range("I12").Value = '//i0.hdslb.com/bfs/openplatform/202403/F86lgbSt1709278264141.jpeg'
$ws4.Range("B13").Value = '''2024-06-15'
$ws4.Range("C13").Value = '苏州·蔚蓝档案ONLY#2024~Game Builders Go!!!!'
$ws4.Range("D13").Value = '城际路21号 苏州汇融广场假日酒店'
$ws4.Range("E13").Value = '2024.06.15 10:00-06.15 17:00'
$ws4.Range("F13").Value = 176
$ws4.Range("G13").Value = 75
$ws4.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=84130'
$ws4.Range("I13").Value = '//i0.hdslb.com/bfs/openplatform/202404/bpTzFcDq1713253785881.jpeg'
$ws4.Range("B14").Value = '''2024-06-16'
$ws4.Range("C14").Value = '苏州·明日方舟ONLY#2024~佑桑柔'
$ws4.Range("D14").Value = '城际路21号 苏州汇融广场假日酒店'
$ws4.Range("E14").Value = '2024.06.16 10:00-06.16 17:00'
$ws4.Range("F14").Value = 147
$ws4.Range("G14").Value = 75
$ws4.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=84046'
$ws4.Range("I14").Value = '//i1.hdslb.com/bfs/openplatform/202405/FxnSAnJm1714978943173.jpeg'
$ws4.Range("B15").Value = '''2024-06-29'
$ws4.Range("C15").Value = '苏州·归离之缘原神only展'
$ws4.Range("D15").Value = '清禾路888号2号楼3楼 格莱美婚礼宴会中心'
$ws4.Range("E15").Value = '2024.06.29 09:30-06.29 18:30'
$ws4.Range("F15").Value = 265
$ws4.Range("G15").Value = 89
$ws4.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=83271'
$ws4.Range("I15").Value = '//i1.hdslb.com/bfs/openplatform/202403/hNkSoRCt1710999968954.png'
$ws4.Range("B16").Value = '''2024-07-06'
$ws4.Range("C16").Value = '苏州·第一届寒假动漫展宅舞比赛-CF01'
$ws4.Range("D16").Value = '润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店'
$ws4.Range("E16").Value = '2024.07.06 10:00-07.06 16:00'
$ws4.Range("F16").Value = 196
$ws4.Range("G16").Value = 49
$ws4.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=80528'
$ws4.Range("I16").Value = '//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg'
$ws4.Range("B17").Value = '''2024-07-20'
$ws4.Range("C17").Value = '苏州·白日梦想7.20全职猎人ONLY展'
$ws4.Range("D17").Value = '金芳路与新发路交叉口东南120米 万龙大厦'
$ws4.Range("E17").Value = '2024.07.20 09:00-07.20 17:00'
$ws4.Range("F17").Value = 394
$ws4.Range("G17").Value = 72
$ws4.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=83508'
$ws4.Range("I17").Value = '//i1.hdslb.com/bfs/openplatform/202403/LC3LtiWw1711633827120.jpeg'
$ws4.Range("B18").Value = '''2024-07-20'
$ws4.Range("C18").Value = '苏州·萤火国潮文化节动漫品牌博览会'
$ws4.Range("D18").Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws4.Range("E18").Value = '2024.07.20 10:00-07.21 17:00'
$ws4.Range("F18").Value = 7170
$ws4.Range("G18").Value = 60
$ws4.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=83301'
$ws4.Range("I18").Value = '//i0.hdslb.com/bfs/openplatform/202403/rV07luU61711274774556.jpeg'
$ws4.Range("B19").Value = '''2024-07-27'
$ws4.Range("C19").Value = '苏州·第一届动漫游戏展'
$ws4.Range("D19").Value = '清禾路886号 尹山湖大剧院'
$ws4.Range("E19").Value = '2024.07.27 10:30-07.27 17:00'
$ws4.Range("F19").Value = 1106
$ws4.Range("G19").Value = 60
$ws4.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=84899'
$ws4.Range("I19").Value = '//i2.hdslb.com/bfs/openplatform/202404/ARz0BVLv1712661597595.jpeg'
$ws4.Range("B20").Value = '''2024-08-03'
$ws4.Range("C20").Value = '苏州·星部落动漫嘉年华'
$ws4.Range("D20").Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Range("E20").Value = '2024.08.03 09:00-08.04 16:00'
$ws4.Range("F20").Value = 676
$ws4.Range("G20").Value = 49
$ws4.Range("H20").Value = 'https://show.bilibili.com/platform/detail.html?id=84858'
$ws4.Range("I20").Value = '//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg'
$ws4.Range("B21").Value = '''2024-08-17'
$ws4.Range("C21").Value = '苏州·ICAN summer World动漫品牌夏游节'
$ws4.Range("D21").Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws4.Range("E21").Value = '2024.08.17 10:00-08.18 17:00'
$ws4.Range("F21").Value = 108
$ws4.Range("G21").Value = 60
$ws4.Range("H21").Value = 'https://show.bilibili.com/platform/detail.html?id=85289'
$ws4.Range("I21").Value = '//i0.hdslb.com/bfs/openplatform/202404/JavlW9fj1714459472747.jpeg'
$ws4.Range("B22").Value = '''2024-10-01'
$ws4.Range("C22").Value = '苏州·第二届Redamancy动漫游戏嘉年华'
$ws4.Range("D22").Value = '长江路436号绿宝广场b1 party king运动街区'
$ws4.Range("E22").Value = '2024.10.01 10:00-10.05 17:00'
$ws4.Range("F22").Value = 58
$ws4.Range("G22").Value = 98
$ws4.Range("H22").Value = 'https://show.bilibili.com/platform/detail.html?id=83576'
$ws4.Range("I22").Value = '//i2.hdslb.com/bfs/openplatform/202403/MKyrtd4c1711689984512.jpeg'
$ws4.Range("B23").Value = '''2024-10-01'
$ws4.Range("C23").Value = '苏州·第十三届理想乡动漫展-同人创作者大会'
$ws4.Range("D23").Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Range("E23").Value = '2024.10.01 10:00-10.03 17:00'
$ws4.Range("F23").Value = 254
$ws4.Range("G23").Value = 39
$ws4.Range("H23").Value = 'https://show.bilibili.com/platform/detail.html?id=83821'
$ws4.Range("I23").Value = '//i0.hdslb.com/bfs/openplatform/202404/OMtuTTFY1711958198579.jpeg'

Write-Output "done"